$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1466.0834
$ws.Range("I15").Value = 1466.0834
$ws.Range("K15").Value = 4398.2502
$ws.Range("M15").Value = -4229.2502
$ws.Range("H17").Value = 1595.6285
$ws.Range("J17").Value = 1829.409
$ws.Range("L17").Value = 5488.227000000001
$ws.Range("N17").Value = -5824.227000000001
$ws.Range("H74").Value = 5473.25
$ws.Range("I74").Value = 5109.7646
$ws.Range("K74").Value = 5109.7646
$ws.Range("M74").Value = -4173.7646
$ws.Range("H77").Value = 5473.25
$ws.Range("I77").Value = 5109.7646
$ws.Range("K77").Value = 25548.823
$ws.Range("M77").Value = -20868.823
$ws.Range("H106").Value = 2248.7646
$ws.Range("I106").Value = 1889.3125
$ws.Range("K106").Value = 1889.3125
$ws.Range("M106").Value = -1258.3125
$ws.Range("H113").Value = 7335.9165
$ws.Range("I113").Value = 5866.6665
$ws.Range("J113").Value = 8805.166999999999
$ws.Range("K113").Value = 5866.6665
$ws.Range("L113").Value = 8805.166999999999
$ws.Range("M113").Value = -2612.6665
$ws.Range("N113").Value = -15313.167
$ws.Range("H138").Value = 6267.8374
$ws.Range("J138").Value = 6764.722
$ws.Range("L138").Value = 20294.166
$ws.Range("N138").Value = -30574.166
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8533555
$ws.Range("I32").Value = 16252340
$ws.Range("K32").Value = 16252340
$ws.Range("M32").Value = -16252053
$ws.Range("H61").Value = 4931.7856
$ws.Range("I61").Value = 3225.8386
$ws.Range("K61").Value = 3225.8386
$ws.Range("M61").Value = -3013.8386
$ws.Range("H74").Value = 2843934
$ws.Range("I74").Value = 3789648.5
$ws.Range("J74").Value = 6790.4546
$ws.Range("K74").Value = 3789648.5
$ws.Range("L74").Value = 6790.4546
$ws.Range("M74").Value = -3788774.5
$ws.Range("N74").Value = -8538.454600000001
$ws.Range("H77").Value = 2843934
$ws.Range("I77").Value = 3789648.5
$ws.Range("J77").Value = 6790.4546
$ws.Range("K77").Value = 18948242.5
$ws.Range("L77").Value = 33952.273
$ws.Range("M77").Value = -18943874.5
$ws.Range("N77").Value = -42688.273
$ws.Range("H82").Value = 20000
$ws.Range("J82").Value = 20000
$ws.Range("L82").Value = 20000
$ws.Range("N82").Value = -20722
$ws.Range("H85").Value = 20000
$ws.Range("J85").Value = 20000
$ws.Range("L85").Value = 20000
$ws.Range("N85").Value = -22496
$ws.Range("H97").Value = 198.75
$ws.Range("I97").Value = 198.75
$ws.Range("K97").Value = 198.75
$ws.Range("M97").Value = 297.25
$ws.Range("H122").Value = 2853.9167
$ws.Range("I122").Value = 2565.1177
$ws.Range("J122").Value = 3555.2856
$ws.Range("K122").Value = 7695.353099999999
$ws.Range("L122").Value = 10665.8568
$ws.Range("M122").Value = -5245.353099999999
$ws.Range("N122").Value = -15565.8568
$ws.Range("H136").Value = 4931.7856
$ws.Range("I136").Value = 3225.8386
$ws.Range("K136").Value = 9677.515800000001
$ws.Range("M136").Value = -7127.515800000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 89999.5
$ws.Range("J21").Value = 89999.5
$ws.Range("L21").Value = 89999.5
$ws.Range("N21").Value = -90471.5
$ws.Range("H82").Value = 35417.582
$ws.Range("J82").Value = 69651
$ws.Range("L82").Value = 69651
$ws.Range("N82").Value = -70417
$ws.Range("H85").Value = 35417.582
$ws.Range("J85").Value = 69651
$ws.Range("L85").Value = 69651
$ws.Range("N85").Value = -72303
$ws.Range("H86").Value = 5000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 5000
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 5000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 25000
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -36232
$ws.Range("H94").Value = 2277.8696
$ws.Range("I94").Value = 845.46155
$ws.Range("J94").Value = 4140
$ws.Range("K94").Value = 845.46155
$ws.Range("L94").Value = 4140
$ws.Range("M94").Value = -394.46155
$ws.Range("N94").Value = -5042
$ws.Range("H97").Value = 14838.143
$ws.Range("I97").Value = 10566
$ws.Range("K97").Value = 10566
$ws.Range("M97").Value = -9575
$ws.Range("H99").Value = 2901.3157
$ws.Range("I99").Value = 1941.0667
$ws.Range("K99").Value = 1941.0667
$ws.Range("M99").Value = -443.0667000000001
$ws.Range("H105").Value = 850.2727
$ws.Range("I105").Value = 850.2727
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 850.2727
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 896.7273
$ws.Range("N105").ClearContents()
$ws.Range("H106").Value = 73335.5
$ws.Range("J106").Value = 73335.5
$ws.Range("L106").Value = 73335.5
$ws.Range("N106").Value = -75859.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16573.055
$ws.Range("J31").Value = 9426.458000000001
$ws.Range("L31").Value = 9426.458000000001
$ws.Range("N31").Value = -10016.458
$ws.Range("H34").Value = 16573.055
$ws.Range("J34").Value = 9426.458000000001
$ws.Range("L34").Value = 9426.458000000001
$ws.Range("N34").Value = -9830.458000000001
$ws.Range("H50").Value = 129799
$ws.Range("J50").Value = 129799
$ws.Range("L50").Value = 129799
$ws.Range("N50").Value = -131049
$ws.Range("H51").Value = 49999.75
$ws.Range("J51").Value = 49999.75
$ws.Range("L51").Value = 49999.75
$ws.Range("N51").Value = -51471.75
$ws.Range("H58").Value = 564126.7
$ws.Range("I58").Value = 774078.9399999999
$ws.Range("K58").Value = 774078.9399999999
$ws.Range("M58").Value = -773875.9399999999
$ws.Range("H59").Value = 129997
$ws.Range("J59").Value = 129997
$ws.Range("L59").Value = 129997
$ws.Range("N59").Value = -132287
$ws.Range("H61").Value = 49999.75
$ws.Range("J61").Value = 49999.75
$ws.Range("L61").Value = 49999.75
$ws.Range("N61").Value = -50695.75
$ws.Range("H74").Value = 44666.61
$ws.Range("J74").Value = 44666.61
$ws.Range("L74").Value = 44666.61
$ws.Range("N74").Value = -46414.61
$ws.Range("H77").Value = 44666.61
$ws.Range("J77").Value = 44666.61
$ws.Range("L77").Value = 133999.83
$ws.Range("N77").Value = -142735.83
$ws.Range("H105").Value = 12111.625
$ws.Range("I105").Value = 12719.134
$ws.Range("J105").Value = 2999
$ws.Range("K105").Value = 12719.134
$ws.Range("L105").Value = 2999
$ws.Range("M105").Value = -10972.134
$ws.Range("N105").Value = -6493
$ws.Range("H134").Value = 345172.5
$ws.Range("I134").Value = 345172.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 1035517.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1032982.5
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 564126.7
$ws.Range("I136").Value = 774078.9399999999
$ws.Range("K136").Value = 2322236.82
$ws.Range("M136").Value = -2319686.82
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4550.3335
$ws.Range("J34").Value = 6500
$ws.Range("L34").Value = 19500
$ws.Range("N34").Value = -19668
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H55").Value = 2449.5
$ws.Range("I55").Value = 2449.5
$ws.Range("K55").Value = 7348.5
$ws.Range("M55").Value = -7171.5
$ws.Range("H113").Value = 14111.429
$ws.Range("J113").Value = 10945.25
$ws.Range("L113").Value = 32835.75
$ws.Range("N113").Value = -37175.75
$ws.Range("H131").Value = 8952
$ws.Range("I131").Value = 1502.5
$ws.Range("J131").Value = 9465.759
$ws.Range("K131").Value = 4507.5
$ws.Range("L131").Value = 28397.277
$ws.Range("M131").Value = 532.5
$ws.Range("N131").Value = -38477.277
$ws.Range("H138").Value = 1779.6875
$ws.Range("I138").Value = 1514.2858
$ws.Range("J138").Value = 3637.5
$ws.Range("K138").Value = 4542.857400000001
$ws.Range("L138").Value = 10912.5
$ws.Range("M138").Value = 597.1425999999992
$ws.Range("N138").Value = -21192.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1958.4889
$ws.Range("I102").Value = 1530.6216
$ws.Range("K102").Value = 1530.6216
$ws.Range("M102").Value = 91.37840000000006
$ws.Range("H105").Value = 85505.86
$ws.Range("J105").Value = 85505.86
$ws.Range("L105").Value = 85505.86
$ws.Range("N105").Value = -92493.86
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H122").Value = 104800.1
$ws.Range("I122").Value = 4999.6665
$ws.Range("J122").Value = 254500.75
$ws.Range("K122").Value = 14998.9995
$ws.Range("L122").Value = 763502.25
$ws.Range("M122").Value = -12548.9995
$ws.Range("N122").Value = -768402.25
$ws.Range("H136").Value = 6080.2
$ws.Range("I136").Value = 5975.75
$ws.Range("K136").Value = 17927.25
$ws.Range("M136").Value = -15377.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H136").Value = 1007418.7
$ws.Range("I136").Value = 2199450.2
$ws.Range("J136").Value = 3602.6316
$ws.Range("K136").Value = 6598350.600000001
$ws.Range("L136").Value = 10807.8948
$ws.Range("M136").Value = -6595800.600000001
